$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / view geometry -------------------------------------------------
# Reposition & resize the workbook window (matches the new bookViews/workbookView
# xWindow/yWindow/windowWidth/windowHeight in workbook.xml).
$win = $wb.Windows.Item(1)
$win.Left   = -120
$win.Top    = -120
$win.Width  = 29040
$win.Height = 15840

# Scroll the sheet back to the top-left corner (A1) so the saved sheetView no
# longer carries a topLeftCell="A16" override.
$win.ScrollRow    = 1
$win.ScrollColumn = 1

# Keep the original selected cell.
$ws.Range("H23").Select()

# --- Column widths -----------------------------------------------------------
# Column A grew wider.
$ws.Columns(1).ColumnWidth = 58.333333333333336

# Columns C:J used to share/overlap a single run of widths; they are now each
# sized individually (a new, much narrower column E appears in the run and the
# others pick up distinct best-fit-like widths).
$ws.Columns(3).ColumnWidth  = 32.666666666666664
$ws.Columns(4).ColumnWidth  = 38.666666666666664
$ws.Columns(5).ColumnWidth  = 7.666666666666667
$ws.Columns(6).ColumnWidth  = 34.833333333333336
$ws.Columns(7).ColumnWidth  = 30.833333333333332
$ws.Columns(8).ColumnWidth  = 30.0
$ws.Columns(9).ColumnWidth  = 34.833333333333336
$ws.Columns(10).ColumnWidth = 33.666666666666664
